$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 5.603999999999999
$ws.Range("C4").Value = -12.829

$ws.Range("C5").Value = -12.413

$ws.Range("B7").Value = 6.957000000000001

$ws.Range("C8").Value = -12.8

$ws.Range("B16").Value = 5.603
$ws.Range("C16").Value = -12.382
